# irrdash3.xlsx — refresh the simulated IRR/cashflow figures on Sheet1
# (columns C "IRR" and D "Mkt cap", rows 4-39) with a newly re-run batch of
# values, and move the active selection to where the author last clicked.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C4").Value  = -121.59310206259966
$ws.Range("D4").Value  = 4611.1790284361814
$ws.Range("C5").Value  = 1392.0699002682704
$ws.Range("D5").Value  = 2567.4289707736116
$ws.Range("C6").Value  = 1838.0876872796198
$ws.Range("D6").Value  = 989.20857123354551
$ws.Range("C7").Value  = 1754.901826303947
$ws.Range("D7").Value  = 221.30914058918347
$ws.Range("C8").Value  = 1710.7683203219235
$ws.Range("D8").Value  = 2106.6256296028787
$ws.Range("C9").Value  = 1647.0656812304687
$ws.Range("D9").Value  = 5294.6326235182451
$ws.Range("C10").Value = 1917.7340557348034
$ws.Range("D10").Value = 4427.5790990897876
$ws.Range("C11").Value = 2029.3531672554727
$ws.Range("D11").Value = 6699.9736305483148
$ws.Range("C12").Value = 1997.8470570346981
$ws.Range("D12").Value = 6197.7326120145926
$ws.Range("C13").Value = 1991.984589130936
$ws.Range("D13").Value = 8111.977184249773
$ws.Range("C14").Value = 2057.1948533202863
$ws.Range("D14").Value = 8099.8737408743837
$ws.Range("C15").Value = 2011.0737052291402
$ws.Range("D15").Value = 8288.8144031864449
$ws.Range("C16").Value = 2007.6284353019623
$ws.Range("D16").Value = 7428.9529537233184
$ws.Range("C17").Value = 2028.547263757755
$ws.Range("D17").Value = 7636.3305582043422
$ws.Range("C18").Value = 2034.06030417581
$ws.Range("D18").Value = 7537.1078262266747
$ws.Range("C19").Value = 1996.6375211769482
$ws.Range("D19").Value = 7239.4228885436014
$ws.Range("C20").Value = 1842.4093558445461
$ws.Range("D20").Value = 7128.8511428279435
$ws.Range("C21").Value = 3128.2932890377224
$ws.Range("D21").Value = 6983.7524039036143
$ws.Range("C22").Value = 1739.2665875421385
$ws.Range("D22").Value = 6736.245310886472
$ws.Range("C23").Value = 1647.2597047795962
$ws.Range("D23").Value = 7188.6301287649376
$ws.Range("C24").Value = -2195.4704725154825
$ws.Range("D24").Value = 7435.8005160484045
$ws.Range("C25").Value = 8383.1529293208077
$ws.Range("D25").Value = 9736.2645221372186
$ws.Range("C26").Value = 416.63938550832677
$ws.Range("D26").Value = 6916.0347688484908
$ws.Range("C27").Value = 357.61370388398888
$ws.Range("D27").Value = 17792.412517860637
$ws.Range("C28").Value = 363.24130408974872
$ws.Range("D28").Value = 21724.322353370055
$ws.Range("C29").Value = 349.11790942905429
$ws.Range("D29").Value = 1858.9630331803621
$ws.Range("C30").Value = 333.88971618004257
$ws.Range("D30").Value = 1984.3939238409062
$ws.Range("C31").Value = 449.03465195398087
$ws.Range("D31").Value = 2044.0408195728162
$ws.Range("C32").Value = 305.07381364064832
$ws.Range("D32").Value = 1996.6431526746217
$ws.Range("C33").Value = 78.665003562786595
$ws.Range("D33").Value = 1898.8909132589124
$ws.Range("C34").Value = -1319.7553754182525
$ws.Range("D34").Value = 1797.543851649908
$ws.Range("D35").Value = -440.11639145949835
$ws.Range("D36").Value = 1515.8089564799636
$ws.Range("D37").Value = 1521.0671101142289
$ws.Range("D38").Value = 1525.3642902194081
$ws.Range("D39").Value = 17770.775179445653

# Move the selection to match where the author left off (T23) on Sheet1.
$ws.Range("T23").Select() | Out-Null
